# ACE_landing_page_data.xlsx update
# - Refresh the ACE KPI table (sheet "ACE_landing_page_data", B2:O7) with new data pull.
# - Refresh the ANSP reference list (sheet "ANSP"): drop "UkSATSE" and add "HCAA"
#   (inserted right after "HASP", before "HungaroControl").

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ACE_landing_page_data")
$ws2 = $wb.Worksheets.Item("ANSP")

# ---------------------------------------------------------------------------
# 1) Update the KPI values on the ACE_landing_page_data sheet
# ---------------------------------------------------------------------------
$ws1.Range("B2").Value2 = 470.2655103
$ws1.Range("C2").Value2 = 8923563036.0
$ws1.Range("D2").Value2 = 18975584.73
$ws1.Range("E2").Value2 = 0.885574066
$ws1.Range("F2").Value2 = 133.2905453
$ws1.Range("G2").Value2 = 319.7523568
$ws1.Range("H2").Value2 = -0.346585581
$ws1.Range("I2").Value2 = 0.034281998
$ws1.Range("J2").Value2 = 0.582888236
$ws1.Range("K2").Value2 = 0.469560726
$ws1.Range("L2").Value2 = -0.023923536
$ws1.Range("M2").Value2 = -0.351540455
$ws1.Range("N2").Value2 = 96.95010124
$ws1.Range("O2").Value2 = 93.29066939

$ws1.Range("B3").Value2 = 719.7048251
$ws1.Range("C3").Value2 = 8627785321.0
$ws1.Range("D3").Value2 = 11987949.81
$ws1.Range("E3").Value2 = 0.602611413
$ws1.Range("F3").Value2 = 136.5574832
$ws1.Range("G3").Value2 = 493.0953049
$ws1.Range("H3").Value2 = -0.252952008
$ws1.Range("I3").Value2 = -0.049996459
$ws1.Range("J3").Value2 = 0.271676721
$ws1.Range("K3").Value2 = 0.251345475
$ws1.Range("L3").Value2 = -0.082098113
$ws1.Range("M3").Value2 = -0.246572025
$ws1.Range("N3").Value2 = 93.73662257
$ws1.Range("O3").Value2 = 58.93699079

$ws1.Range("B4").Value2 = 963.3983799
$ws1.Range("C4").Value2 = 9081845434.0
$ws1.Range("D4").Value2 = 9426884.685
$ws1.Range("E4").Value2 = 0.481570777
$ws1.Range("F4").Value2 = 148.7713285
$ws1.Range("G4").Value2 = 654.4690687
$ws1.Range("H4").Value2 = 1.213430842
$ws1.Range("I4").Value2 = -0.041892622
$ws1.Range("J4").Value2 = -0.56713923
$ws1.Range("K4").Value2 = -0.50710518
$ws1.Range("L4").Value2 = 0.084551656
$ws1.Range("M4").Value2 = 1.219649314
$ws1.Range("N4").Value2 = 98.66976125
$ws1.Range("O4").Value2 = 46.34589106

$ws1.Range("B5").Value2 = 435.2511773
$ws1.Range("C5").Value2 = 9478943216.0
$ws1.Range("D5").Value2 = 21778098.97
$ws1.Range("E5").Value2 = 0.977025437
$ws1.Range("F5").Value2 = 137.1731145
$ws1.Range("G5").Value2 = 294.8524637
$ws1.Range("H5").Value2 = -0.002038008
$ws1.Range("I5").Value2 = 0.014498579
$ws1.Range("J5").Value2 = 0.016570358
$ws1.Range("K5").Value2 = 0.010387739
$ws1.Range("L5").Value2 = 0.006381333
$ws1.Range("M5").Value2 = -0.001117714
$ws1.Range("N5").Value2 = 102.9840324
$ws1.Range("O5").Value2 = 107.0688182

$ws1.Range("B6").Value2 = 436.1400342
$ws1.Range("C6").Value2 = 9343476090.0
$ws1.Range("D6").Value2 = 21423110.37
$ws1.Range("E6").Value2 = 0.966980694
$ws1.Range("F6").Value2 = 136.3033177
$ws1.Range("G6").Value2 = 295.182393
$ws1.Range("H6").Value2 = -0.036186772
$ws1.Range("I6").Value2 = 0.01512249
$ws1.Range("J6").Value2 = 0.05323569
$ws1.Range("K6").Value2 = 0.049273702
$ws1.Range("L6").Value2 = -0.00017276
$ws1.Range("M6").Value2 = -0.030874655
$ws1.Range("N6").Value2 = 101.512249
$ws1.Range("O6").Value2 = 105.323569

$ws1.Range("B7").Value2 = 452.5150948
$ws1.Range("C7").Value2 = 9204284392.0
$ws1.Range("D7").Value2 = 20340281.46
$ws1.Range("E7").Value2 = 0.921571457
$ws1.Range("F7").Value2 = 136.3268695
$ws1.Range("G7").Value2 = 304.5863927
$ws1.Range("H7").Value2 = -0.03586958
$ws1.Range("I7").Value2 = 0.008963882
$ws1.Range("J7").Value2 = 0.046501449
$ws1.Range("K7").Value2 = 0.043322939
$ws1.Range("L7").Value2 = 0.012469209
$ws1.Range("M7").Value2 = -0.038898468
$ws1.Range("N7").Value2 = 100.0
$ws1.Range("O7").Value2 = 100.0

# ---------------------------------------------------------------------------
# 2) Update the ANSP reference list: insert "HCAA" right after "HASP" (row 18),
#    pushing the remaining ANSPs down one row, then drop the old trailing
#    "UkSATSE" entry that used to be in row 40.
# ---------------------------------------------------------------------------
$ws2.Rows.Item(19).Insert()
$ws2.Cells.Item(18, 1).Copy()
$ws2.Cells.Item(19, 1).PasteSpecial(-4122) # xlPasteFormats
$ws2.Application.CutCopyMode = $false
$ws2.Cells.Item(19, 1).Value2 = "HCAA"

# Row 40 now carries the distinctive "last row" formatting that used to sit on
# the real last row (now shifted to 41, "UkSATSE"). Move that formatting back
# onto the new last row (40, "SMATSA") before deleting the now-duplicate
# "UkSATSE" row.
$ws2.Cells.Item(41, 1).Copy()
$ws2.Cells.Item(40, 1).PasteSpecial(-4122) # xlPasteFormats
$ws2.Application.CutCopyMode = $false

$ws2.Rows.Item(41).Delete()

# ---------------------------------------------------------------------------
# 3) Restore the selections / active sheet shown when the file was saved.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A2:A40").Select()

$ws1.Activate()
$ws1.Range("A1:O7").Select()
